# Update the "register" sheet's sample row (row 2) with new test-user data,
# re-point the existing hyperlink + add hyperlinks for the two password
# columns, and switch the active sheet/selection from "login" to "register".

$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("login")
$register = $wb.Worksheets.Item("register")

# --- Update the sample data on the "register" sheet ---------------------
$register.Range("A2").Value = "Peter"
$register.Range("B2").Value = "White"
$register.Range("C2").Value = "peterwh@gmail.com"
$register.Range("D2").Value = "Peter@26"
$register.Range("E2").Value = "Peter@26"

# --- Re-point the existing E-mail hyperlink to the new address ----------
$register.Range("C2").Hyperlinks.Delete()
$register.Hyperlinks.Add($register.Range("C2"), "mailto:peterwh@gmail.com")
$register.Range("C2").Style = "Hyperlink"

# --- Add hyperlinks on the Password / C-Password cells -------------------
$register.Hyperlinks.Add($register.Range("D2"), "mailto:Peter@26")
$register.Range("D2").Style = "Hyperlink"

$register.Hyperlinks.Add($register.Range("E2"), "mailto:Peter@26")
$register.Range("E2").Style = "Hyperlink"

# --- Make "register" the active sheet/tab and select F2 -----------------
$register.Activate()
$register.Range("F2").Select()
